$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 12 (the 키움/코셈 entry dated 2024-02-13) - remaining rows shift up.
$ws.Rows.Item(12).Delete()
